$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.96219999999999
$ws.Range("A12").Value = -22.7487
$ws.Range("E12").Value = 12.77919999999999
$ws.Range("E14").Value = 13.79420000000001
$ws.Range("E22").Value = 11.9332
